$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.247.74'
$ws.Range('E2').Value = '  +2.94%  '
$ws.Range('D3').Value = '1.716.95'
$ws.Range('E3').Value = '  +3.31%  '
$ws.Range('D4').Value = "'1.000"
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = "'238.71"
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('D6').Value = "'1.001"
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('D7').Value = "'0.4724"
$ws.Range('E7').Value = '  -1.24%  '
$ws.Range('D8').Value = "'0.2614"
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = "'0.06191"
$ws.Range('E9').Value = '  +0.30%  '
$ws.Range('D10').Value = '1.715.16'
$ws.Range('E10').Value = '  +3.20%  '
$ws.Range('D11').Value = "'0.07079"
$ws.Range('D12').Value = "'15.22"
$ws.Range('E12').Value = '  +3.09%  '
$ws.Range('D13').Value = "'0.5900"
$ws.Range('E13').Value = '  +0.50%  '
$ws.Range('D14').Value = "'4.407"
$ws.Range('E14').Value = '  +1.18%  '
$ws.Range('D15').Value = "'75.93"
$ws.Range('E15').Value = '  +1.88%  '
$ws.Range('D16').Value = "'1.001"
$ws.Range('E16').Value = '  +0.10%  '
$ws.Range('B17').Value = 'BinanceUSD'
$ws.Range('C17').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D17').Value = "'1.001"
$ws.Range('E17').Value = '  +0.12%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '26.251.79'
$ws.Range('E18').Value = '  +2.95%  '
$ws.Range('D19').Value = "'0.000006807"
$ws.Range('E19').Value = '  +0.61%  '
$ws.Range('D20').Value = "'11.49"
$ws.Range('E20').Value = '  +0.69%  '
$ws.Range('D21').Value = '1.935.10'
$ws.Range('E21').Value = '  +3.59%  '
$ws.Range('D22').Value = "'4.544"
$ws.Range('E22').Value = '  +2.54%  '
$ws.Range('D23').Value = "'8.706"
$ws.Range('E23').Value = '  +0.88%  '
$ws.Range('D24').Value = "'5.278"
$ws.Range('E24').Value = '  +0.20%  '
$ws.Range('D25').Value = "'134.95"
$ws.Range('E25').Value = '  +1.57%  '
$ws.Range('D26').Value = "'15.16"
$ws.Range('E26').Value = '  +1.05%  '
$ws.Range('D27').Value = "'1.403"
$ws.Range('E27').Value = '  +0.85%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = "'1.754"
$ws.Range('E28').Value = '  +3.44%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = "'107.26"
$ws.Range('E29').Value = '  +2.46%  '
$ws.Range('D30').Value = "'3.968"
$ws.Range('E30').Value = '  +0.45%  '
$ws.Range('D31').Value = "'3.670"
$ws.Range('E31').Value = '  +0.32%  '
$ws.Range('D32').Value = "'0.07745"
$ws.Range('E32').Value = '  +1.55%  '
$ws.Range('D33').Value = "'0.04421"
$ws.Range('E33').Value = '  +4.84%  '
$ws.Range('D34').Value = "'2.617"
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('D35').Value = "'0.9705"
$ws.Range('E35').Value = '  +2.11%  '
$ws.Range('D36').Value = "'0.6138"
$ws.Range('E36').Value = '  +0.72%  '
$ws.Range('D37').Value = "'0.9228"
$ws.Range('E37').Value = '  +7.92%  '
$ws.Range('D38').Value = "'112.24"
$ws.Range('E38').Value = '  +15.68%  '
$ws.Range('D39').Value = "'2.449"
$ws.Range('E39').Value = '  -5.82%  '
$ws.Range('D40').Value = "'1.912"
$ws.Range('E40').Value = '  +3.05%  '
$ws.Range('E41').Value = '  +0.09%  '
$ws.Range('D42').Value = "'0.01468"
$ws.Range('E42').Value = '  -0.59%  '
$ws.Range('D43').Value = "'5.363"
$ws.Range('E43').Value = '  +13.51%  '
$ws.Range('D44').Value = "'0.3804"
$ws.Range('E44').Value = '  +1.32%  '
$ws.Range('D45').Value = "'0.1162"
$ws.Range('E45').Value = '  +3.95%  '
$ws.Range('D46').Value = "'6.241"
$ws.Range('E46').Value = '  +0.59%  '
$ws.Range('D47').Value = "'0.05262"
$ws.Range('E47').Value = '  +0.18%  '
$ws.Range('D48').Value = "'30.15"
$ws.Range('E48').Value = '  +2.05%  '
$ws.Range('D49').Value = "'7.699"
$ws.Range('E49').Value = '  +5.45%  '
$ws.Range('D50').Value = "'0.3379"
$ws.Range('E50').Value = '  +1.59%  '
$ws.Range('E51').Value = '  +0.11%  '
